$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("English", 20.71308487690542),
    @("Chinese", 20.5303461260024),
    @("Spanish", 6.179306121893809),
    @("Arabic", 4.17247083327666),
    @("German", 3.913821941223753),
    @("Japanese", 3.630077487128431),
    @("Malay-Indonesian", 3.288115946473386),
    @("Russian", 2.934420383339788),
    @("Portuguese", 2.719508630757407),
    @("French", 2.433508848536508),
    @("Turkish", 2.059928609444937),
    @("Italian", 1.792389085418489),
    @("Korean", 1.667105117990741),
    @("Dutch", 1.165014191415372),
    @("Polish", 0.9722930816214356),
    @("Persian", 0.9608754149229302),
    @("Thai", 0.9175711547039331),
    @("Urdu", 0.9120956405487437),
    @("Bengali", 0.8803059684204186),
    @("Vietnamese", 0.8627947775661128)
)

# Overwrite the existing data rows (2 through 21) with the new sorted values,
# preserving the existing cell formatting (style stays attached to the cell).
$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# Rows 22 and 23 no longer exist in the table, delete their entire rows so the
# sheet dimension shrinks back down to A1:B21 and no leftover formatting remains.
$ws.Range("A22:B23").EntireRow.Delete()
